# "line of balance" chart (plotly/jupyter) rework:
# the "space" sheet's name column becomes a plain numeric space id (1, 2)
# instead of floor_1/floor_2 text, and the "design" sheet's space column
# is repointed at that same numeric id, with quantity collapsed to a
# constant 8 (the per-space-unit productivity denominator used downstream).

$wb = $excel.ActiveWorkbook

# --- "space" sheet -------------------------------------------------------
$wsSpace = $wb.Worksheets.Item("space")
$wsSpace.Range("A2").Value = 1
$wsSpace.Range("A3").Value = 2

# --- "design" sheet -------------------------------------------------------
$wsDesign = $wb.Worksheets.Item("design")
$wsDesign.Range("C1").Value = "space"

$wsDesign.Range("B2").Value = 1
$wsDesign.Range("C2").Value = 8

$wsDesign.Range("B3").Value = 1
$wsDesign.Range("C3").Value = 8

$wsDesign.Range("B4").Value = 1
$wsDesign.Range("C4").Value = 8

$wsDesign.Range("B5").Value = 1
$wsDesign.Range("C5").Value = 8

$wsDesign.Range("B6").Value = 2
$wsDesign.Range("C6").Value = 8

$wsDesign.Range("B7").Value = 2
$wsDesign.Range("C7").Value = 8

# --- selection / active-sheet bookkeeping ---------------------------------
# "space" no longer the tab in focus; leave its cursor on the row below the
# data instead of the old "whole column" selection.
[void]$wsSpace.Select()
[void]$wsSpace.Range("A4").Select()

# "design" becomes the active tab, cursor parked on the (now numeric)
# space column.
[void]$wsDesign.Select()
[void]$wsDesign.Range("C2:C7").Select()
